$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44034
$ws.Range("C2").Value = 29588
$ws.Range("D2").Value = 764
$ws.Range("E2").Value = 3491
$ws.Range("H2").Value = 0.12
$ws.Range("B4").Value = 44034
$ws.Range("C4").Value = '''219128'
$ws.Range("D4").Value = '''18803'
$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43
$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217
$ws.Range("B7").Value = 44034
$ws.Range("C7").Value = 84417
$ws.Range("D7").Value = 888
$ws.Range("E7").Value = 16107
$ws.Range("F7").Value = 319
$ws.Range("G7").Value = 19.08
$ws.Range("H7").Value = 35.92
$ws.Range("B8").Value = 44034
$ws.Range("C8").Value = '''35578'
$ws.Range("D8").Value = '''260'
$ws.Range("E8").Value = '''886'
$ws.Range("F8").Value = '''5'
$ws.Range("H8").Value = 1.92
$ws.Range("B9").Value = 44034
$ws.Range("C9").Value = 24540
$ws.Range("D9").Value = 677
$ws.Range("E9").Value = 2044
$ws.Range("G9").Value = 12.43
$ws.Range("H9").Value = 14.69
$ws.Range("K9").Value = 16446
$ws.Range("L9").Value = 632
$ws.Range("B13").Value = 44034
$ws.Range("C13").Value = 17828
$ws.Range("D13").Value = 591
$ws.Range("E13").Value = 322
$ws.Range("B16").Value = 44034
$ws.Range("C16").Value = 70413
$ws.Range("D16").Value = 1325
$ws.Range("E16").Value = 20886
$ws.Range("F16").Value = 561
$ws.Range("G16").Value = 43.68
$ws.Range("H16").Value = 44.38
$ws.Range("K16").Value = 47812
$ws.Range("L16").Value = 1264
$ws.Range("B19").Value = 44033
$ws.Range("C19").Value = 47071
$ws.Range("D19").Value = 1423
$ws.Range("E19").Value = 20775
$ws.Range("F19").Value = 717
$ws.Range("G19").Value = 44.14
$ws.Range("H19").Value = 50.39
$ws.Range("B20").Value = 44034
$ws.Range("C20").Value = 150609
$ws.Range("D20").Value = 2974
$ws.Range("E20").Value = 3246
$ws.Range("F20").Value = 86
$ws.Range("G20").Value = 4.31
$ws.Range("H20").Value = 3.47
$ws.Range("K20").Value = 75236
$ws.Range("L20").Value = 2479
$ws.Range("C21").Value = 100483
$ws.Range("D21").Value = 7063
$ws.Range("E21").Value = 14038
$ws.Range("G21").Value = 29.99
$ws.Range("K21").Value = 46806
$ws.Range("B23").Value = 44034
$ws.Range("C23").Value = 2813
$ws.Range("D23").Value = 42
$ws.Range("G23").Value = 0.53
$ws.Range("B24").Value = 44034
$ws.Range("G24").Value = 11.39
$ws.Range("K24").Value = 1326
$ws.Range("B25").Value = 44033
$ws.Range("C25").Value = 46203
$ws.Range("E25").Value = 6375
$ws.Range("F25").Value = 651
$ws.Range("G25").Value = 13.8
$ws.Range("H25").Value = 18.46
$ws.Range("B26").Value = 44034
$ws.Range("C26").Value = 24104
$ws.Range("D26").Value = 308
$ws.Range("E26").Value = 1780
$ws.Range("F26").Value = 64
$ws.Range("G26").Value = 9.029999999999999
$ws.Range("H26").Value = 21.33
$ws.Range("K26").Value = 19703
$ws.Range("L26").Value = 300
$ws.Range("B27").Value = 44034
$ws.Range("C27").Value = 41698
$ws.Range("D27").Value = 1771
$ws.Range("E27").Value = 2027
$ws.Range("G27").Value = 6.17
$ws.Range("H27").Value = 6.89
$ws.Range("K27").Value = 32858
$ws.Range("L27").Value = 1713
$ws.Range("B29").Value = 44034
$ws.Range("C29").Value = 75171
$ws.Range("D29").Value = 6056
$ws.Range("E29").Value = 21592
$ws.Range("F29").Value = 2415
$ws.Range("G29").Value = 28.72
$ws.Range("H29").Value = 39.88
$ws.Range("B30").Value = 44034
$ws.Range("C30").Value = 99354
$ws.Range("D30").Value = 3558
$ws.Range("E30").Value = 36693
$ws.Range("F30").Value = 1803
$ws.Range("G30").Value = 45.36
$ws.Range("H30").Value = 51.05
$ws.Range("K30").Value = 80885
$ws.Range("L30").Value = 3532
$ws.Range("B31").Value = 44033
$ws.Range("C31").Value = 413576
$ws.Range("D31").Value = 7870
$ws.Range("E31").Value = 11396
$ws.Range("F31").Value = 660
$ws.Range("G31").Value = 4.31
$ws.Range("H31").Value = 8.6
$ws.Range("K31").Value = 264380
$ws.Range("L31").Value = 7672
$ws.Range("B32").Value = 44034
$ws.Range("C32").Value = 58673
$ws.Range("D32").Value = 2666
$ws.Range("E32").Value = 6696
$ws.Range("F32").Value = 377
$ws.Range("G32").Value = 11.41
$ws.Range("B33").Value = 44034
$ws.Range("C33").Value = 2132
$ws.Range("D33").Value = 19
$ws.Range("E33").Value = 49
$ws.Range("G33").Value = 3.7
$ws.Range("K33").Value = 1323
$ws.Range("B34").Value = 44034
$ws.Range("C34").Value = 44847
$ws.Range("D34").Value = 865
$ws.Range("E34").Value = 7066
$ws.Range("F34").Value = 197
$ws.Range("G34").Value = 17.43
$ws.Range("H34").Value = 23.15
$ws.Range("K34").Value = 40535
$ws.Range("L34").Value = 851
$ws.Range("B35").Value = 44034
$ws.Range("C35").Value = 152302
$ws.Range("D35").Value = 3335
$ws.Range("E35").Value = 39406
$ws.Range("F35").Value = 1525
$ws.Range("G35").Value = 25.87
$ws.Range("H35").Value = 45.73
$ws.Range("B37").Value = 44034
$ws.Range("C37").Value = 6295
$ws.Range("D37").Value = 402
$ws.Range("E37").Value = 325
$ws.Range("G37").Value = 6.01
$ws.Range("H37").Value = 2.25
$ws.Range("K37").Value = 5404
$ws.Range("L37").Value = 400
$ws.Range("C41").Value = 40000
$ws.Range("D41").Value = 811
$ws.Range("E41").Value = 3280
$ws.Range("H41").Value = 4.69
$ws.Range("B42").Value = 44034
$ws.Range("C42").Value = 105001
$ws.Range("D42").Value = 1698
$ws.Range("E42").Value = 17314
$ws.Range("F42").Value = 536
$ws.Range("G42").Value = 23.99
$ws.Range("H42").Value = 32.66
$ws.Range("K42").Value = 72174
$ws.Range("L42").Value = 1641
$ws.Range("B43").Value = 44034
$ws.Range("C43").Value = 165301
$ws.Range("D43").Value = 7347
$ws.Range("E43").Value = 27704
$ws.Range("F43").Value = 2024
$ws.Range("G43").Value = 16.76
$ws.Range("B44:H44").Value = "'"
$ws.Range("J44").Value = $false
$ws.Range("O44").Value = 'An error occurred. ... TimeoutException('''', None, None)'
$ws.Range("B45").Value = 44034
$ws.Range("C45").Value = 47961
$ws.Range("D45").Value = 1552
$ws.Range("E45").Value = 9842
$ws.Range("F45").Value = 152
$ws.Range("G45").Value = 20.52
$ws.Range("H45").Value = 9.789999999999999
$ws.Range("B46").Value = 44034
$ws.Range("C46").Value = 114320
$ws.Range("D46").Value = 8468
$ws.Range("E46").Value = 10731
$ws.Range("H46").Value = 8.18
$ws.Range("B47").Value = 44034
$ws.Range("C47").Value = 16911
$ws.Range("D47").Value = 370
$ws.Range("E47").Value = 4771
$ws.Range("F47").Value = 147
$ws.Range("G47").Value = 30.42
$ws.Range("H47").Value = 39.73
$ws.Range("K47").Value = 15682
$ws.Range("L47").Value = 370
$ws.Range("B48").Value = 44034
$ws.Range("C48").Value = 36063
$ws.Range("D48").Value = 1159
$ws.Range("E48").Value = 8118
$ws.Range("F48").Value = 389
$ws.Range("G48").Value = 31.68
$ws.Range("H48").Value = 36.22
$ws.Range("K48").Value = 25629
$ws.Range("L48").Value = 1074

Write-Output "Applied all cell updates"
